$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.522.01'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '2.265.37'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'119.23"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.83%  '
$ws.Range('D6').Value = "'265.08"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  +2.99%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('D10').Value = "'47.41"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('D11').Value = "'0.0942"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = "'9.11"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.10%  '
$ws.Range('E13').Value = '  -1.53%  '
$ws.Range('D14').Value = "'15.32"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('D15').Value = "'0.905"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.84%  '
$ws.Range('D16').Value = '2.604.63'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '2.269.60'
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = '43.498.77'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = "'0.0000110"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').Value = "'6.86"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('D21').Value = "'72.03"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = "'235.21"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('D24').Value = "'9.50"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.48%  '
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('E26').Value = '  +1.95%  '
$ws.Range('D27').Value = "'11.86"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.93%  '
$ws.Range('D28').Value = "'41.48"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('E30').Value = '  +0.13%  '
$ws.Range('D31').Value = "'171.82"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('D32').Value = "'21.71"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('D33').Value = "'0.0912"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').Value = "'5.68"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = "'4.32"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +15.69%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = "'0.130"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.31%  '
$ws.Range('D37').Value = "'0.0376"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.99%  '
$ws.Range('D38').Value = "'4.55"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('E39').Value = '  +1.21%  '
$ws.Range('D40').Value = "'2.54"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.05%  '
$ws.Range('D41').Value = "'13.92"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.85%  '
$ws.Range('D42').Value = "'74.19"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('E46').Value = '  -8.18%  '
$ws.Range('D47').Value = "'74.11"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +41.59%  '
$ws.Range('D48').Value = "'8.51"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.84%  '
$ws.Range('D49').Value = "'1.26"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('D51').Value = "'101.28"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.56%  '
